$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unmerge existing merged cells before restructuring
$ws.Range("A5:A6").UnMerge()
$ws.Range("A3:A4").UnMerge()

# Insert a new row before row 5, pushing old row5 (proceedings-article) down to row6,
# and old row6 down to row7. This gives room for the new "tex-math" row.
$ws.Rows.Item(5).Insert()

# Row 3: journal-article:None | MathML | 13
$ws.Cells.Item(3, 2).Value = "MathML"
$ws.Cells.Item(3, 3).Value = 13

# Row 4: (blank A) | No | 185
$ws.Cells.Item(4, 2).Value = "No"
$ws.Cells.Item(4, 3).Value = 185

# Row 5 (new): (blank A) | tex-math | 2
$ws.Cells.Item(4, 1).Copy($ws.Cells.Item(5, 1))
$ws.Cells.Item(4, 2).Copy($ws.Cells.Item(5, 2))
$ws.Cells.Item(5, 2).Value = "tex-math"
$ws.Cells.Item(5, 3).Value = 2

# Row 6: proceedings-article:None | HTML | 1
$ws.Cells.Item(6, 2).Value = "HTML"
$ws.Cells.Item(6, 3).Value = 1

# Row 7: (blank A) | No | 49
$ws.Cells.Item(7, 2).Value = "No"
$ws.Cells.Item(7, 3).Value = 49

# Re-create merges for the new layout
$ws.Range("A3:A5").Merge()
$ws.Range("A6:A7").Merge()
